# Update the crypto price/volume snapshot data (GitHub Actions refresh).
# Cells D (Price) and E (Volume(1h)) are plain text values (e.g. "291.12",
# "-3.29%"), so each new value is written with a leading apostrophe to force
# Excel to store it as literal text rather than re-interpreting it as a
# number/percentage. The style is reset to "Normal" right after so the cell
# keeps its original (unstyled) appearance instead of picking up a "Text"
# number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$updates = @{
    "D2"  = "291.12";  "E2"  = "-3.29%"
    "D3"  = "30.67";   "E3"  = "-6.42%"
    "D4"  = "4.945";   "E4"  = "0.16%"
    "D5"  = "0.07207"; "E5"  = "-6.44%"
    "D6"  = "1.821";   "E6"  = "-7.89%"
    "D7"  = "7.692";   "E7"  = "-1.90%"
    "D8"  = "3.765";   "E8"  = "-0.85%"
    "D9"  = "0.8963";  "E9"  = "-2.54%"
    "D10" = "0.1650";  "E10" = "-5.90%"
    "D11" = "0.07725"; "E11" = "-0.35%"
    "D12" = "0.07940"; "E12" = "-7.91%"
    "D13" = "0.03037"; "E13" = "-4.36%"
    "E14" = "-0.19%"
    "E15" = "-0.52%"
    "D16" = "0.005693"; "E16" = "-3.73%"
    "D18" = "3.473";    "E18" = "0.43%"
    "E19" = "-3.35%"
    "E20" = "-0.95%"
    "D21" = "0.1313";   "E21" = "-0.94%"
    "D22" = "4.034";    "E22" = "-6.11%"
    "D23" = "0.2388";   "E23" = "19.95%"
    "D24" = "0.04493";  "E24" = "-0.37%"
    "D25" = "0.001215"
    "D26" = "0.004008"; "E26" = "-9.30%"
    "E27" = "-0.01%"
    "D39" = "0.01590";        "E39" = "-6.78%"
    "D40" = "0.04402";        "E40" = "-6.24%"
    "D41" = "0.007242";       "E41" = "-3.23%"
    "D42" = "0.01002"
    "D43" = "0.1309";         "E43" = "-3.23%"
    "D44" = "0.002051";       "E44" = "-12.02%"
    "D45" = "0.009511";       "E45" = "-8.97%"
    "D46" = "0.00005925";     "E46" = "-5.19%"
    "D47" = "0.00000000750";  "E47" = "0.00%"
    "E48" = "172.74%"
    "D49" = "0.003002";       "E49" = "-3.33%"
    "D50" = "0.00002101";     "E50" = "0.00%"
    "D51" = "0.0002001";      "E51" = "0.00%"
}

foreach ($cellRef in $updates.Keys) {
    Set-TextValue $cellRef $updates[$cellRef]
}
